$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.954.40"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.359.46"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'0.662"
$ws.Range("E5").Value = "  +2.87%  "

$ws.Range("D6").Value = "'235.46"
$ws.Range("E6").Value = "  +1.24%  "

$ws.Range("D7").Value = "'72.59"
$ws.Range("E7").Value = "  +10.20%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  +23.59%  "

$ws.Range("D10").Value = "'0.0988"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").Value = "'28.03"
$ws.Range("E11").Value = "  +4.54%  "

$ws.Range("D12").Value = "2.715.10"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "'16.87"
$ws.Range("E14").Value = "  +10.43%  "

$ws.Range("D15").Value = "'6.69"
$ws.Range("E15").Value = "  +9.34%  "

$ws.Range("D16").Value = "'0.883"
$ws.Range("E16").Value = "  +5.54%  "

$ws.Range("D17").Value = "2.330.11"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "43.935.46"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("D20").Value = "'76.30"
$ws.Range("E20").Value = "  +3.79%  "

$ws.Range("D21").Value = "'6.33"
$ws.Range("E21").Value = "  +2.27%  "

$ws.Range("D22").Value = "'251.57"
$ws.Range("E22").Value = "  +1.07%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").Value = "'3.80"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("D26").Value = "'10.50"
$ws.Range("E26").Value = "  +6.66%  "

$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").Value = "'22.50"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").Value = "'173.10"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("E30").Value = "  +9.27%  "

$ws.Range("D31").Value = "'0.131"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  +4.69%  "

$ws.Range("D33").Value = "'5.20"
$ws.Range("E33").Value = "  +4.47%  "

$ws.Range("D34").Value = "'0.0711"
$ws.Range("E34").Value = "  +4.09%  "

$ws.Range("E35").Value = "  +4.08%  "

$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("D37").Value = "'2.43"
$ws.Range("E37").Value = "  +2.19%  "

$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").Value = "'0.0274"
$ws.Range("E39").Value = "  +8.76%  "

$ws.Range("D40").Value = "'19.26"
$ws.Range("E40").Value = "  +11.28%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.94"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("D44").Value = "'0.0971"
$ws.Range("E44").Value = "  +2.15%  "

$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("E46").Value = "  +13.78%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'97.95"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'4.43"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("D49").Value = "1.438.09"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  +0.81%  "

$ws.Range("D51").Value = "2.583.67"
$ws.Range("E51").Value = "  +1.66%  "
